$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "57.678.91"
$ws.Range("E2").Value = "  -2.49%  "

$ws.Range("D3").Value = "2.442.35"
$ws.Range("E3").Value = "  -3.72%  "

$ws.Range("E4").Value = "  +0.02%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "520.52"
$ws.Range("E5").Value = "  -1.16%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "128.04"
$ws.Range("E6").Value = "  -5.15%  "

$ws.Range("E7").Value = "  -0.01%  "

$ws.Range("E8").Value = "  -1.41%  "

$ws.Range("E9").Value = "  -2.12%  "

$ws.Range("E10").Value = "  -2.27%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "4.92"
$ws.Range("E11").Value = "  -5.69%  "

$ws.Range("E12").Value = "  -4.73%  "

$ws.Range("D13").Value = "2.875.47"
$ws.Range("E13").Value = "  -3.79%  "

$ws.Range("D14").Value = "57.617.57"
$ws.Range("E14").Value = "  -2.48%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "21.47"
$ws.Range("E15").Value = "  -3.98%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000131"
$ws.Range("E16").Value = "  -3.77%  "

$ws.Range("D17").Value = "2.443.13"
$ws.Range("E17").Value = "  -3.67%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "10.35"
$ws.Range("E18").Value = "  -3.87%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.08"
$ws.Range("E19").Value = "  -2.98%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "310.88"
$ws.Range("E20").Value = "  -4.06%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.07"
$ws.Range("E21").Value = "  -0.45%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.998"
$ws.Range("E22").Value = "  -0.12%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "64.71"
$ws.Range("E23").Value = "  -0.86%  "

$ws.Range("B24").Value = "Binance-PegBSC-USD"
$ws.Range("C24").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.00"
$ws.Range("E24").Value = "  +0.37%  "

$ws.Range("B25").Value = "Polygon"
$ws.Range("C25").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.399"
$ws.Range("E25").Value = "  -2.81%  "

$ws.Range("B26").Value = "WrappedeETH"
$ws.Range("C26").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Range("D26").Value = "2.560.17"
$ws.Range("E26").Value = "  -3.34%  "

$ws.Range("E27").Value = "  -3.73%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.19"
$ws.Range("E28").Value = "  -4.36%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "173.15"
$ws.Range("E29").Value = "  +1.88%  "

$ws.Range("E30").Value = "  -3.90%  "

$ws.Range("E31").Value = "  -3.42%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.13"

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.12"
$ws.Range("E33").Value = "  -9.09%  "

$ws.Range("E34").Value = "  +0.01%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.999"
$ws.Range("E35").Value = "  -0.09%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "17.74"
$ws.Range("E36").Value = "  -3.04%  "

$ws.Range("E37").Value = "  -7.47%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.73"
$ws.Range("E38").Value = "  -6.49%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "36.27"
$ws.Range("E39").Value = "  -1.27%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.792"
$ws.Range("E40").Value = "  +0.77%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.42"
$ws.Range("E41").Value = "  -6.13%  "

$ws.Range("E42").Value = "  -3.85%  "

$ws.Range("E43").Value = "  -3.86%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "4.73"
$ws.Range("E44").Value = "  -7.07%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0916"
$ws.Range("E45").Value = "  -0.84%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "249.91"
$ws.Range("E46").Value = "  -10.87%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "119.21"
$ws.Range("E47").Value = "  -11.64%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0489"
$ws.Range("E48").Value = "  -3.58%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0209"
$ws.Range("E49").Value = "  -3.91%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "16.89"
$ws.Range("E50").Value = "  -5.47%  "

$ws.Range("E51").Value = "  -0.65%  "
